$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "Qui sint nemo in vol "
$ws.Range("C8").Value = "A fugiat aute nesci "

$ws.Range("B11").Value = "Yo, Et eiusmod odit ad q certifico haber recibido el carnet de identificación como Servidor de la Empresa y a la vez me responsabilizo por el buen uso del mismo y en caso de la finalización de la relación laboral a la devolución del mismo. `n`nLibre y voluntariamente me comprometo a depositar en la cuenta de la Empresa Eléctrica Regional Centro Sur C.A.  el valor de 10usd, correspodiente a la reposición por pérdida del carnet de identificación que he recibido. "

$ws.Range("B17").Value = "23 de octubre de 2020 "
$ws.Range("C17").Value = "Consequatur numquam  "
$ws.Range("D17").Value = "Quidem quidem placea "
$ws.Range("E17").Value = "Ea quas non sit haru "

$ws.Range("D22").Value = "Et vel est nulla aut "

$ws.Range("D23").Value = "23 de octubre de 2020 "
